# Update F-column "want to go" counts on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1030
$ws1.Range("F3").Value = 2158
$ws1.Range("F4").Value = 8
$ws1.Range("F5").Value = 474

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1030
$ws4.Range("F5").Value = 2158
$ws4.Range("F6").Value = 8
$ws4.Range("F7").Value = 474
